$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two remaining "catatan" lines with the corrected wording.
$ws.Range("A6").Value = "1. Sesuaikan kolom peserta (A) berdasarkan sasaran : - penduduk = nik, - keluarga = no. kk, - rumah tangga = no. rtm, - kelompok = kode"
$ws.Range("A7").Value = "2. Kolom Peserta (A) dan kolom NIK (C) wajib di isi, yang lain jika kosong data diambil dari data penduduk berdasarkan kolom NIK ©"

# Match the number format already used for the A4/A5 footnote cells
# (text format, so leading digits like "1." and "2." aren't mangled).
$ws.Range("A6:A7").NumberFormat = "@"

# Remove the now-obsolete rows 8 and 9 (old points 3 and 4).
$ws.Range("A8:A9").EntireRow.Delete()

# Leave the selection where the user ended up after the edit.
$ws.Range("A8").Select()
